$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N9").Value = 9519.58
$ws.Range("O9").Value = 9278.5

$ws.Range("N12").Value = 417357.09
$ws.Range("O12").Value = 357010.17

$ws.Range("O28").Value = 62355
